# Implemented Controller actions and page for Deposit and updated estimations.
#
# This updates the "Estimacija" tracking sheet:
#   - Adds the "Real time (minutes)" actuals for the two newly finished rows
#     (row 19 "Dodavanje rute za Deposit sredstava u WalletController" and
#      row 20 "Dodavanje stranice za Deposit sredstava u MVC aplikaciju")
#   - Updates the sheet view's scroll position / selection to reflect where
#     the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the newly-tracked "real time" actuals -------------------------
$ws.Range("D19").Value = 5
$ws.Range("D20").Value = 90

# --- Update view state: scroll so row 4 is at the top, select B23 ----------
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1

$ws.Range("B23").Select() | Out-Null
